# Update workbook values in sheets "展览" and "全部类型" to reflect the
# newly generated output data (counts of interested attendees / min prices).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 14052
$ws1.Range("G2").Value = 70
$ws1.Range("G3").Value = 60
$ws1.Range("F6").Value = 1217
$ws1.Range("F7").Value = 1047
$ws1.Range("F8").Value = 13923
$ws1.Range("F9").Value = 14942
$ws1.Range("F10").Value = 4
$ws1.Range("F11").Value = 13
$ws1.Range("F25").Value = 5813
$ws1.Range("F28").Value = 5442
$ws1.Range("F32").Value = 320

# ---- Sheet "全部类型" (all types) ----
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 14052
$ws4.Range("G2").Value = 70
$ws4.Range("G3").Value = 60
$ws4.Range("F7").Value = 1217
$ws4.Range("F8").Value = 1047
$ws4.Range("F9").Value = 13923
$ws4.Range("F10").Value = 14942
$ws4.Range("F11").Value = 4
$ws4.Range("F12").Value = 13
$ws4.Range("F27").Value = 5813
$ws4.Range("F30").Value = 5442
$ws4.Range("F34").Value = 320
